$p = $ppt.ActivePresentation

function Find-ShapeByText($slide, $needle) {
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -like "*$needle*") {
                return $shp
            }
        }
    }
    return $null
}

# --- Story "E-6" slide: add an "Overlaps with C-6" note --------------------
# The Notes box on this slide ends with a paragraph holding just a single
# placeholder space character; append the cross-reference note to it.
# Find the slide whose title run reads exactly "Story E-6" (avoid any other
# slide that happens to mention the code elsewhere).
$sE6 = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $titleShape = Find-ShapeByText $slide "Story "
    if ($titleShape -ne $null -and $titleShape.TextFrame.TextRange.Text -eq "Story E-6") {
        $sE6 = $slide
        break
    }
}

$notesShapeE6 = Find-ShapeByText $sE6 "Notes"
$notesRangeE6 = $notesShapeE6.TextFrame.TextRange
$lastParaE6 = $notesRangeE6.Paragraphs($notesRangeE6.Paragraphs().Count)
$lastParaE6.Text = $lastParaE6.Text + "Overlaps with C-6"

# --- Story "C-6" slide: rename header + append a sign-off note -------------
$sC6 = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $titleShape = Find-ShapeByText $slide "Story "
    if ($titleShape -ne $null -and $titleShape.TextFrame.TextRange.Text -eq "Story C-6") {
        $sC6 = $slide
        break
    }
}

$headerShapeC6 = Find-ShapeByText $sC6 "Staff Account"
$headerShapeC6.TextFrame.TextRange.Text = "Account Management"

$notesShapeC6 = Find-ShapeByText $sC6 "Story overlaps with E-6"
$notesRangeC6 = $notesShapeC6.TextFrame.TextRange
$lastParaC6 = $notesRangeC6.Paragraphs($notesRangeC6.Paragraphs().Count)
$lastParaC6.InsertAfter("`rSigned off by team 118 (Kane N8866856) on Oct 5 2016") | Out-Null
